$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.479.85'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.674.73'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '648.72'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.01'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.06'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.438'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000230'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D13").Value = '4.289.21'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.45'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '3.673.34'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").Value = '69.433.55'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.98'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.43'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.53'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.71'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.641'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.35'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").Value = '3.815.41'
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +1.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.80'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.91'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.61'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.66'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.79%  '
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.99'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.46'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.57'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").Value = '3.660.40'
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("E36").Value = '  +2.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.35'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.91'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '179.16'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0892'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.930'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.77'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.70'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.26'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.99'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.62%  '
$ws.Range("E49").Value = '  -4.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.79'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("E51").Value = '  -4.22%  '
